# Update view-count figures (column F) across the four sheets to match
# the newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 8831
$ws1.Range("F4").Value  = 1949
$ws1.Range("F7").Value  = 2115
$ws1.Range("F9").Value  = 66
$ws1.Range("F11").Value = 66
$ws1.Range("F14").Value = 74
$ws1.Range("F16").Value = 8726
$ws1.Range("F19").Value = 200
$ws1.Range("F26").Value = 66
$ws1.Range("F30").Value = 25
$ws1.Range("F34").Value = 2199
$ws1.Range("F35").Value = 867
$ws1.Range("F38").Value = 6
$ws1.Range("F40").Value = 241
$ws1.Range("F41").Value = 171
$ws1.Range("F43").Value = 459
$ws1.Range("F44").Value = 78
$ws1.Range("F45").Value = 98
$ws1.Range("F47").Value = 3987

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 403

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 322

# --- Sheet 4: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 403
$ws4.Range("F5").Value  = 8831
$ws4.Range("F7").Value  = 322
$ws4.Range("F8").Value  = 1949
$ws4.Range("F11").Value = 2115
$ws4.Range("F17").Value = 66
$ws4.Range("F18").Value = 74
$ws4.Range("F20").Value = 8726
$ws4.Range("F22").Value = 200
$ws4.Range("F28").Value = 66
$ws4.Range("F34").Value = 2199
$ws4.Range("F35").Value = 867
$ws4.Range("F38").Value = 6
$ws4.Range("F40").Value = 241
$ws4.Range("F42").Value = 171
$ws4.Range("F44").Value = 98
$ws4.Range("F45").Value = 3987
